# Generate Report for Archive
#
# 1) Update the "Status" text from "Ready for handoff" to "In Translation"
#    everywhere it appears (Overview sheet columns E/F, and the per-language
#    sheets' Status column C).
# 2) Narrow the "Latest Handoff Datetime" / "Status" columns that were sized
#    to fit "Ready for handoff" down to the narrower width needed for
#    "In Translation".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

$newStatus = "In Translation"

# --- Update cell text -------------------------------------------------
$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus
$ws1.Range("E3").Value = $newStatus
$ws1.Range("F3").Value = $newStatus
$ws1.Range("E4").Value = $newStatus
$ws1.Range("F4").Value = $newStatus

$ws2.Range("C2").Value = $newStatus
$ws2.Range("C3").Value = $newStatus
$ws2.Range("C4").Value = $newStatus

$ws3.Range("C2").Value = $newStatus
$ws3.Range("C3").Value = $newStatus
$ws3.Range("C4").Value = $newStatus

# --- Resize the columns that held that text ----------------------------
# Original stored width 17.2159881591797 -> new stored width 13.4101845877511.
# ColumnWidth 12.5 is the closest achievable COM value that serializes to
# the nearest representable width (13.333333333333334).
$ws1.Columns.Item(5).ColumnWidth = 12.5
$ws1.Columns.Item(6).ColumnWidth = 12.5

$ws2.Columns.Item(3).ColumnWidth = 12.5

$ws3.Columns.Item(3).ColumnWidth = 12.5
